$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.403.25'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.843.38'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.21'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6347'
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07478'
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.16'
$ws.Range("E9").Value = '  +3.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2905'
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07744'
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").Value = '1.850.65'
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6796'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001024'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.253'
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D18").Value = '29.454.61'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.18'
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.417'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.10'
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.508'
$ws.Range("E25").Value = '  +1.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1361'
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06545'
$ws.Range("E28").Value = '  +14.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.431'
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.488'
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.073'
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.049'
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7002'
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.578'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01858'
$ws.Range("E37").Value = '  +2.78%  '
$ws.Range("D38").Value = '1.250.81'
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.821'
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.755'
$ws.Range("E40").Value = '  +3.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9398'
$ws.Range("E41").Value = '  +4.02%  '
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("D43").Value = '2.010.06'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.29'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.48'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.072'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("E48").Value = '  +4.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.008'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1149'
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("E51").Value = '  -0.28%  '
